# Horarios actualizados Linea 141 - 363
# Refreshes the scraped schedule data across all three sheets
# (LP1912, LP1912-215, 6203-6173) with the new scrape taken at 03:55:13.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

$newTime = "03:55:13"

# ---------------------------------------------------------------
# Sheet 1: LP1912  (7 data rows, rows 6-12)
# ---------------------------------------------------------------
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 7"

$sheet1Rows = @(
    @($newTime, "04:01", "81_EL PELIGRO", 6,   "LP1912"),
    @($newTime, "04:34", "215_ALUAR",     39,  "LP1912"),
    @($newTime, "04:53", "11_ETCHEVERRY", 58,  "LP1912"),
    @($newTime, "05:16", "17_ROMERO",     81,  "LP1912"),
    @($newTime, "05:22", "23_HERNANDEZ",  87,  "LP1912"),
    @($newTime, "05:35", "215B_EL PATO",  100, "LP1912"),
    @($newTime, "05:46", "15_ABASTO",     111, "LP1912")
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet 2: LP1912-215  (2 data rows, rows 6-7)
# ---------------------------------------------------------------
$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$sheet2Rows = @(
    @($newTime, "04:34", "215_ALUAR",    39,  "LP1912"),
    @($newTime, "05:35", "215B_EL PATO", 100, "LP1912")
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet 3: 6203-6173  (1 data row, row 6) - previously had no data
# rows, so the header row (row 5) needs to be created first by
# copying the formatting used on the other sheets' header rows.
# ---------------------------------------------------------------
$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 1"

$ws1.Range("A5:E5").Copy()
$ws3.Range("A5:E5").PasteSpecial(-4122)   # xlPasteFormats

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"

$ws3.Cells.Item(6, 1).Value = $newTime
$ws3.Cells.Item(6, 2).Value = "05:44"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 109
$ws3.Cells.Item(6, 5).Value = "L6173"

Write-Host "Schedule refresh complete."
